$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.10419689114334
$ws.Range("C2").Value = 1.983950157387443
$ws.Range("D2").Value = -19.10419689114334
$ws.Range("E2").Value = -19.10419689114334
$ws.Range("F2").Value = -19.10419689114334
$ws.Range("G2").Value = -19.10419689114334
$ws.Range("H2").Value = -19.10419689114334
$ws.Range("I2").Value = -19.10419689114334
$ws.Range("J2").Value = -19.10419689114334
$ws.Range("K2").Value = -19.10419689114334
$ws.Range("B3").Value = -19.10419689114334
$ws.Range("C3").Value = -19.10419689114334
$ws.Range("D3").Value = -19.10419689114334
$ws.Range("E3").Value = -19.10419689114334
$ws.Range("F3").Value = -19.10419689114334
$ws.Range("G3").Value = -19.10419689114334
$ws.Range("H3").Value = -19.10419689114334
$ws.Range("I3").Value = 1.269408545074324
$ws.Range("J3").Value = -19.10419689114334
$ws.Range("K3").Value = -19.10419689114334
$ws.Range("B4").Value = -19.10419689114334
$ws.Range("C4").Value = 1.969981818933036
$ws.Range("D4").Value = 1.657644041953099
$ws.Range("E4").Value = -19.10419689114334
$ws.Range("F4").Value = 3.43210274158827
$ws.Range("G4").Value = -19.10419689114334
$ws.Range("H4").Value = 1.485232071945687
$ws.Range("I4").Value = -19.10419689114334
$ws.Range("J4").Value = -19.10419689114334
$ws.Range("K4").Value = -19.10419689114334
$ws.Range("B5").Value = -19.10419689114334
$ws.Range("C5").Value = 1.665912858815526
$ws.Range("D5").Value = -19.10419689114334
$ws.Range("E5").Value = -19.10419689114334
$ws.Range("F5").Value = -19.10419689114334
$ws.Range("G5").Value = 2.904977356507852
$ws.Range("H5").Value = -19.10419689114334
$ws.Range("I5").Value = -19.10419689114334
$ws.Range("J5").Value = -19.10419689114334
$ws.Range("K5").Value = -19.10419689114334
$ws.Range("B6").Value = -19.10419689114334
$ws.Range("C6").Value = -19.10419689114334
$ws.Range("D6").Value = -19.10419689114334
$ws.Range("E6").Value = -19.10419689114334
$ws.Range("F6").Value = -19.10419689114334
$ws.Range("G6").Value = -19.10419689114334
$ws.Range("H6").Value = -19.10419689114334
$ws.Range("I6").Value = -19.10419689114334
$ws.Range("J6").Value = -19.10419689114334
$ws.Range("K6").Value = -19.10419689114334
$ws.Range("B7").Value = 2.433880709363736
$ws.Range("C7").Value = -19.10419689114334
$ws.Range("D7").Value = -19.10419689114334
$ws.Range("E7").Value = -19.10419689114334
$ws.Range("F7").Value = -19.10419689114334
$ws.Range("G7").Value = -19.10419689114334
$ws.Range("H7").Value = -19.10419689114334
$ws.Range("I7").Value = -19.10419689114334
$ws.Range("J7").Value = -19.10419689114334
$ws.Range("K7").Value = -19.10419689114334
$ws.Range("B8").Value = -19.10419689114334
$ws.Range("C8").Value = -19.10419689114334
$ws.Range("D8").Value = -19.10419689114334
$ws.Range("E8").Value = 1.798638208088598
$ws.Range("F8").Value = -19.10419689114334
$ws.Range("G8").Value = -19.10419689114334
$ws.Range("H8").Value = -19.10419689114334
$ws.Range("I8").Value = -19.10419689114334
$ws.Range("J8").Value = -19.10419689114334
$ws.Range("K8").Value = -19.10419689114334
$ws.Range("B9").Value = 3.867552402424593
$ws.Range("C9").Value = -19.10419689114334
$ws.Range("D9").Value = -19.10419689114334
$ws.Range("E9").Value = -19.10419689114334
$ws.Range("F9").Value = -19.10419689114334
$ws.Range("G9").Value = -19.10419689114334
$ws.Range("H9").Value = -19.10419689114334
$ws.Range("I9").Value = -19.10419689114334
$ws.Range("J9").Value = -19.10419689114334
$ws.Range("K9").Value = -19.10419689114334
$ws.Range("B10").Value = -19.10419689114334
$ws.Range("C10").Value = -19.10419689114334
$ws.Range("D10").Value = -19.10419689114334
$ws.Range("E10").Value = -19.10419689114334
$ws.Range("F10").Value = -19.10419689114334
$ws.Range("G10").Value = -19.10419689114334
$ws.Range("H10").Value = -19.10419689114334
$ws.Range("I10").Value = 1.827646485411568
$ws.Range("J10").Value = -19.10419689114334
$ws.Range("K10").Value = 2.23652959176523
$ws.Range("B11").Value = -19.10419689114334
$ws.Range("C11").Value = -19.10419689114334
$ws.Range("D11").Value = -19.10419689114334
$ws.Range("E11").Value = 2.909359213827974
$ws.Range("F11").Value = -19.10419689114334
$ws.Range("G11").Value = 2.809582586690374
$ws.Range("H11").Value = -19.10419689114334
$ws.Range("I11").Value = -19.10419689114334
$ws.Range("J11").Value = -19.10419689114334
$ws.Range("K11").Value = 1.912987366698852
$ws.Range("B12").Value = -19.10419689114334
$ws.Range("C12").Value = -19.10419689114334
$ws.Range("D12").Value = -19.10419689114334
$ws.Range("E12").Value = -19.10419689114334
$ws.Range("F12").Value = -19.10419689114334
$ws.Range("G12").Value = -19.10419689114334
$ws.Range("H12").Value = -19.10419689114334
$ws.Range("I12").Value = -19.10419689114334
$ws.Range("J12").Value = -19.10419689114334
$ws.Range("K12").Value = -19.10419689114334
$ws.Range("B13").Value = -19.10419689114334
$ws.Range("C13").Value = -19.10419689114334
$ws.Range("D13").Value = -19.10419689114334
$ws.Range("E13").Value = 2.538170839852489
$ws.Range("F13").Value = -19.10419689114334
$ws.Range("G13").Value = -19.10419689114334
$ws.Range("H13").Value = -19.10419689114334
$ws.Range("I13").Value = -19.10419689114334
$ws.Range("J13").Value = -19.10419689114334
$ws.Range("K13").Value = 1.704724089983745
$ws.Range("B14").Value = -19.10419689114334
$ws.Range("C14").Value = -19.10419689114334
$ws.Range("D14").Value = 1.518084657890762
$ws.Range("E14").Value = -19.10419689114334
$ws.Range("F14").Value = -19.10419689114334
$ws.Range("G14").Value = -19.10419689114334
$ws.Range("H14").Value = -19.10419689114334
$ws.Range("I14").Value = -19.10419689114334
$ws.Range("J14").Value = -19.10419689114334
$ws.Range("K14").Value = 1.877936758784984
$ws.Range("B15").Value = -19.10419689114334
$ws.Range("C15").Value = -19.10419689114334
$ws.Range("D15").Value = 1.783771115853536
$ws.Range("E15").Value = -19.10419689114334
$ws.Range("F15").Value = -19.10419689114334
$ws.Range("G15").Value = -19.10419689114334
$ws.Range("H15").Value = -19.10419689114334
$ws.Range("I15").Value = -19.10419689114334
$ws.Range("J15").Value = -19.10419689114334
$ws.Range("K15").Value = -19.10419689114334
$ws.Range("B16").Value = -19.10419689114334
$ws.Range("C16").Value = -19.10419689114334
$ws.Range("D16").Value = -19.10419689114334
$ws.Range("E16").Value = -19.10419689114334
$ws.Range("F16").Value = -19.10419689114334
$ws.Range("G16").Value = -19.10419689114334
$ws.Range("H16").Value = -19.10419689114334
$ws.Range("I16").Value = -19.10419689114334
$ws.Range("J16").Value = 4.321925662895021
$ws.Range("K16").Value = -19.10419689114334
$ws.Range("B17").Value = -19.10419689114334
$ws.Range("C17").Value = 2.156702548302345
$ws.Range("D17").Value = 1.870996854789465
$ws.Range("E17").Value = -19.10419689114334
$ws.Range("F17").Value = -19.10419689114334
$ws.Range("G17").Value = -19.10419689114334
$ws.Range("H17").Value = 1.911790448521852
$ws.Range("I17").Value = 2.154657215951135
$ws.Range("J17").Value = -19.10419689114334
$ws.Range("K17").Value = -19.10419689114334
$ws.Range("B18").Value = -19.10419689114334
$ws.Range("C18").Value = -19.10419689114334
$ws.Range("D18").Value = -19.10419689114334
$ws.Range("E18").Value = -19.10419689114334
$ws.Range("F18").Value = -19.10419689114334
$ws.Range("G18").Value = -19.10419689114334
$ws.Range("H18").Value = 2.031801030328082
$ws.Range("I18").Value = 2.03208898204092
$ws.Range("J18").Value = -19.10419689114334
$ws.Range("K18").Value = -19.10419689114334
$ws.Range("B19").Value = -19.10419689114334
$ws.Range("C19").Value = -19.10419689114334
$ws.Range("D19").Value = 2.026281293309503
$ws.Range("E19").Value = -19.10419689114334
$ws.Range("F19").Value = -19.10419689114334
$ws.Range("G19").Value = -19.10419689114334
$ws.Range("H19").Value = 1.724132545816459
$ws.Range("I19").Value = 1.941259502201348
$ws.Range("J19").Value = -19.10419689114334
$ws.Range("K19").Value = -19.10419689114334
$ws.Range("B20").Value = -19.10419689114334
$ws.Range("C20").Value = 1.060843799474255
$ws.Range("D20").Value = 1.488629656546191
$ws.Range("E20").Value = -19.10419689114334
$ws.Range("F20").Value = 3.202633858350728
$ws.Range("G20").Value = -19.10419689114334
$ws.Range("H20").Value = 1.65043702260965
$ws.Range("I20").Value = 0.728373308466025
$ws.Range("J20").Value = -19.10419689114334
$ws.Range("K20").Value = 2.197335585998641
$ws.Range("B21").Value = -19.10419689114334
$ws.Range("C21").Value = 1.268048595425511
$ws.Range("D21").Value = -19.10419689114334
$ws.Range("E21").Value = 1.677931898878275
$ws.Range("F21").Value = -19.10419689114334
$ws.Range("G21").Value = 2.459193519941169
$ws.Range("H21").Value = 1.538708870920713
$ws.Range("I21").Value = -19.10419689114334
$ws.Range("J21").Value = -19.10419689114334
$ws.Range("K21").Value = -19.10419689114334
